$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "60.709.96"
$ws.Range("E2").Value2 = "  -1.76%  "

$ws.Range("D3").Value2 = "3.383.91"
$ws.Range("E3").Value2 = "  -2.09%  "

$ws.Range("E4").Value2 = "  +0.00%  "

$ws.Range("D5").Value2 = "'569.59"
$ws.Range("E5").Value2 = "  -2.19%  "

$ws.Range("D6").Value2 = "'141.63"
$ws.Range("E6").Value2 = "  -3.61%  "

$ws.Range("E7").Value2 = "  +0.06%  "

$ws.Range("D8").Value2 = "3.383.39"
$ws.Range("E8").Value2 = "  -2.12%  "

$ws.Range("E9").Value2 = "  -0.19%  "

$ws.Range("E10").Value2 = "  -2.31%  "

$ws.Range("E11").Value2 = "  -2.19%  "

$ws.Range("D12").Value2 = "'0.395"
$ws.Range("E12").Value2 = "  +0.86%  "

$ws.Range("D13").Value2 = "3.964.78"
$ws.Range("E13").Value2 = "  -2.01%  "

$ws.Range("D14").Value2 = "'28.32"
$ws.Range("E14").Value2 = "  +1.53%  "

$ws.Range("E15").Value2 = "  +1.52%  "

$ws.Range("D17").Value2 = "3.385.58"

$ws.Range("D18").Value2 = "60.804.89"
$ws.Range("E18").Value2 = "  -1.78%  "

$ws.Range("D19").Value2 = "'6.25"
$ws.Range("E19").Value2 = "  -0.08%  "

$ws.Range("D20").Value2 = "'14.02"
$ws.Range("E20").Value2 = "  -2.68%  "

$ws.Range("D21").Value2 = "'9.00"
$ws.Range("E21").Value2 = "  -5.74%  "

$ws.Range("D22").Value2 = "'385.43"
$ws.Range("E22").Value2 = "  -1.00%  "

$ws.Range("D23").Value2 = "'0.560"
$ws.Range("E23").Value2 = "  -1.04%  "

$ws.Range("D24").Value2 = "'73.71"
$ws.Range("E24").Value2 = "  +0.12%  "

$ws.Range("D25").Value2 = "'0.998"
$ws.Range("E25").Value2 = "  -0.10%  "

$ws.Range("E26").Value2 = "  -5.49%  "

$ws.Range("D27").Value2 = "3.523.59"
$ws.Range("E27").Value2 = "  -2.08%  "

$ws.Range("E28").Value2 = "  -1.68%  "

$ws.Range("E29").Value2 = "  +0.51%  "

$ws.Range("D30").Value2 = "'7.39"
$ws.Range("E30").Value2 = "  -4.46%  "

$ws.Range("E31").Value2 = "  -3.10%  "

$ws.Range("E32").Value2 = "  -2.06%  "

$ws.Range("D33").Value2 = "'1.43"
$ws.Range("E33").Value2 = "  -4.22%  "

$ws.Range("E35").Value2 = "  -2.09%  "

$ws.Range("D37").Value2 = "'166.62"
$ws.Range("E37").Value2 = "  -0.08%  "

$ws.Range("D38").Value2 = "3.414.49"
$ws.Range("E38").Value2 = "  -2.04%  "

$ws.Range("D39").Value2 = "'5.00"
$ws.Range("E39").Value2 = "  -3.16%  "

$ws.Range("E40").Value2 = "  -5.26%  "

$ws.Range("D41").Value2 = "'27.91"
$ws.Range("E41").Value2 = "  +2.28%  "

$ws.Range("E42").Value2 = "  -1.15%  "

$ws.Range("E43").Value2 = "  -3.12%  "

$ws.Range("E44").Value2 = "  +0.01%  "

$ws.Range("E45").Value2 = "  -1.63%  "

$ws.Range("E47").Value2 = "  -3.14%  "

$ws.Range("D48").Value2 = "2.543.76"
$ws.Range("E48").Value2 = "  -1.16%  "

$ws.Range("E49").Value2 = "  -3.82%  "

$ws.Range("D50").Value2 = "'23.54"
$ws.Range("E50").Value2 = "  +2.05%  "

$ws.Range("E51").Value2 = "  -1.20%  "
